$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (shifts existing C:J to D:K),
# carrying over column B's formatting/width like Excel's UI "Insert" does.
$ws.Columns("C:C").Insert()

# The new column should look like its left neighbour (header fill/border style)
$ws.Range("D1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0
$ws.Columns("C:C").ColumnWidth = $ws.Columns("B:B").ColumnWidth

# New "SHIPMENT TYPE" column content
$ws.Range("C1").Value = "SHIPMENT TYPE"
$ws.Range("C2").Value = "DOMESTIC"
$ws.Range("C3").Value = "INTERNATIONAL - EXPORT"

# Reset the view: no frozen/scrolled topLeftCell, selection moves to D7
$ws.Range("D7").Select()
